# Particle contact base resolution algos
# Adds a new diary entry (row 20) to Sheet1, matching the style of the
# previous entry (row 19), and nudges a couple of cosmetic view/format
# properties to line up with the author's saved session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New diary entry -------------------------------------------------
$ws.Range("A20").Value = "31 loka"
$ws.Range("B20").Value = "12.00-13.00, 16.45-"
$ws.Range("C20").Value = "Tsempring, rajoitteet kytketyissä kappaleissa"

# Match formatting used by the rest of the table:
#  - column B ("Kello") is a wrapped, time-formatted cell (style index 3)
#  - column C ("Oppimisen sisältö") is a wrapped general-text cell (style index 2)
$ws.Range("B20").NumberFormat = "h:mm"
$ws.Range("B20").WrapText = $true
$ws.Range("C20").WrapText = $true

# New row renders as a two-line wrapped row, same height as row 19.
$ws.Range("A20:G20").RowHeight = 29

# The last cell of the previous row (F19, "META") also picks up the
# shared wrap-text style used throughout the table.
$ws.Range("F19").WrapText = $true

# --- Minor column width touch-ups (sub-pixel nudges from the session) -
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 40.666666666666664
$ws.Columns.Item(5).ColumnWidth = 33.666666666666664
$ws.Columns.Item(6).ColumnWidth = 31.666666666666668

# --- Selection moved to the newly-entered row -------------------------
$ws.Range("D20").Select() | Out-Null
